$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.544.64'
$ws.Range('E2').Value = '  +5.17%  '
$ws.Range('D3').Value = '3.634.72'
$ws.Range('E3').Value = '  +5.47%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '592.71'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.97%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '191.84'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('E7').Value = '  +2.00%  '
$ws.Range('D8').Value = '3.626.22'
$ws.Range('E8').Value = '  +5.33%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +2.52%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.667'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +3.02%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '58.34'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +3.29%  '
$ws.Range('E13').Value = '  +3.84%  '
$ws.Range('E14').Value = '  +4.81%  '
$ws.Range('D15').Value = '4.214.88'
$ws.Range('E15').Value = '  +5.60%  '
$ws.Range('E16').Value = '  +5.64%  '
$ws.Range('D17').Value = '3.631.34'
$ws.Range('E17').Value = '  +5.50%  '
$ws.Range('D18').Value = '70.512.66'
$ws.Range('E18').Value = '  +5.17%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.70'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +4.94%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  +4.20%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '488.89'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.65%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '19.47'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +15.29%  '
$ws.Range('E24').Value = '  +2.46%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.46'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.24%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '91.01'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +6.21%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range('E29').Value = '  +5.58%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '33.15'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +5.58%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.83'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +9.38%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '633.80'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +7.27%  '
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('E34').Value = '  +7.18%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '66.15'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '38.81'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +6.27%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.412'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +6.64%  '
$ws.Range('D38').Value = '0.0₃0824'
$ws.Range('E38').Value = '  +6.72%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -1.07%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '3.58'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').Value = '3.310.56'
$ws.Range('E42').Value = '  +3.48%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.10'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +6.31%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +10.17%  '
$ws.Range('E45').Value = '  +4.89%  '
$ws.Range('E46').Value = '  +2.68%  '
$ws.Range('E47').Value = '  +2.45%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.12'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +4.36%  '
$ws.Range('E49').Value = '  -2.27%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '3.30'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +3.81%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '142.24'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.78%  '
